$wb = $excel.ActiveWorkbook

# ======================= Sheet1 (Overview) =======================
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Cells.Item(5, 1).Value = "e91463e1-37c1-4132-936d-4741263aa477.md"
$ws1.Hyperlinks.Add($ws1.Cells.Item(5, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cd726b8b47065f6a56789dc7e1d5036ee5c2f35/e2e/e91463e1-37c1-4132-936d-4741263aa477.md", "", "", "e2e\e91463e1-37c1-4132-936d-4741263aa477.md")
$ws1.Cells.Item(5, 3).Value = ".md"
$ws1.Cells.Item(5, 4).Value = ""
$ws1.Cells.Item(5, 5).Value = "Ready for handoff"
$ws1.Cells.Item(5, 6).Value = "Ready for handoff"
$ws1.Cells.Item(5, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws1.Cells.Item(5, 7).Value = "2016-10-19 12:32:14"

$ws1.Cells.Item(6, 1).Value = "72a22f36-255e-4335-94e1-dc1ccefda6fd.png"
$ws1.Hyperlinks.Add($ws1.Cells.Item(6, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cd726b8b47065f6a56789dc7e1d5036ee5c2f35/e2e/72a22f36-255e-4335-94e1-dc1ccefda6fd.png", "", "", "e2e\72a22f36-255e-4335-94e1-dc1ccefda6fd.png")
$ws1.Cells.Item(6, 3).Value = ".png"
$ws1.Cells.Item(6, 4).Value = ""
$ws1.Cells.Item(6, 5).Value = "Ready for handoff"
$ws1.Cells.Item(6, 6).Value = "Ready for handoff"
$ws1.Cells.Item(6, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws1.Cells.Item(6, 7).Value = "2016-10-19 12:32:14"

$ws1.Cells.Item(7, 1).Value = "86614a44-b657-4e47-b440-569090e5a0bf.png"
$ws1.Hyperlinks.Add($ws1.Cells.Item(7, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cd726b8b47065f6a56789dc7e1d5036ee5c2f35/e2e/86614a44-b657-4e47-b440-569090e5a0bf.png", "", "", "e2e\86614a44-b657-4e47-b440-569090e5a0bf.png")
$ws1.Cells.Item(7, 3).Value = ".png"
$ws1.Cells.Item(7, 4).Value = ""
$ws1.Cells.Item(7, 5).Value = "Ready for handoff"
$ws1.Cells.Item(7, 6).Value = "Ready for handoff"
$ws1.Cells.Item(7, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws1.Cells.Item(7, 7).Value = "2016-10-19 12:32:14"

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G7"))

# ======================= zh-cn =======================
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Add($ws2.Cells.Item(5, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3855458c6c3725275fbb4796c7832864e5cc411c/e2e/e91463e1-37c1-4132-936d-4741263aa477.md", "", "", "e91463e1-37c1-4132-936d-4741263aa477.md")
$ws2.Cells.Item(5, 2).Value = ".md"
$ws2.Cells.Item(5, 3).Value = "Ready for handoff"
$ws2.Cells.Item(5, 4).Value = "e2e"
$ws2.Cells.Item(5, 5).Value = "ht"
$ws2.Cells.Item(5, 6).Value = "False"
$ws2.Cells.Item(5, 7).Value = "e91463e1-37c1-4132-936d-4741263aa477.eb0d211ad163a74bdec9c4fa15f72fa75ab483ad.zh-cn.xlf"
$ws2.Cells.Item(5, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(5, 8).Value = "2016-10-19 12:32:00"
$ws2.Cells.Item(5, 9).Value = ""
$ws2.Cells.Item(5, 10).Value = ""
$ws2.Cells.Item(5, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(5, 11).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(5, 12).Value = ""
$ws2.Cells.Item(5, 13).Value = "True"
$ws2.Cells.Item(5, 14).Value = ""
$ws2.Cells.Item(5, 15).Value = "False"
$ws2.Cells.Item(5, 16).Value = ""

$ws2.Hyperlinks.Add($ws2.Cells.Item(6, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3855458c6c3725275fbb4796c7832864e5cc411c/e2e/72a22f36-255e-4335-94e1-dc1ccefda6fd.png", "", "", "72a22f36-255e-4335-94e1-dc1ccefda6fd.png")
$ws2.Cells.Item(6, 2).Value = ".png"
$ws2.Cells.Item(6, 3).Value = "Ready for handoff"
$ws2.Cells.Item(6, 4).Value = "e2e"
$ws2.Cells.Item(6, 5).Value = "ht"
$ws2.Cells.Item(6, 6).Value = "False"
$ws2.Cells.Item(6, 7).Value = "e80e7f845ffd185e6ca385a6da93685d70922b36.png"
$ws2.Cells.Item(6, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(6, 8).Value = "2016-10-19 12:32:00"
$ws2.Cells.Item(6, 9).Value = ""
$ws2.Cells.Item(6, 10).Value = ""
$ws2.Cells.Item(6, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(6, 11).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(6, 12).Value = ""
$ws2.Cells.Item(6, 13).Value = "True(Dependency)"
$ws2.Cells.Item(6, 14).Value = "e2e\e91463e1-37c1-4132-936d-4741263aa477.md"
$ws2.Cells.Item(6, 15).Value = "False"
$ws2.Cells.Item(6, 16).Value = ""

$ws2.Hyperlinks.Add($ws2.Cells.Item(7, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3855458c6c3725275fbb4796c7832864e5cc411c/e2e/86614a44-b657-4e47-b440-569090e5a0bf.png", "", "", "86614a44-b657-4e47-b440-569090e5a0bf.png")
$ws2.Cells.Item(7, 2).Value = ".png"
$ws2.Cells.Item(7, 3).Value = "Ready for handoff"
$ws2.Cells.Item(7, 4).Value = "e2e"
$ws2.Cells.Item(7, 5).Value = "ht"
$ws2.Cells.Item(7, 6).Value = "False"
$ws2.Cells.Item(7, 7).Value = "3c4c472bb964b0033e821dd8fb8d4d40e13f94cd.png"
$ws2.Cells.Item(7, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(7, 8).Value = "2016-10-19 12:32:00"
$ws2.Cells.Item(7, 9).Value = ""
$ws2.Cells.Item(7, 10).Value = ""
$ws2.Cells.Item(7, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(7, 11).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(7, 12).Value = ""
$ws2.Cells.Item(7, 13).Value = "True(Dependency)"
$ws2.Cells.Item(7, 14).Value = "e2e\e91463e1-37c1-4132-936d-4741263aa477.md"
$ws2.Cells.Item(7, 15).Value = "False"
$ws2.Cells.Item(7, 16).Value = ""

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P7"))
$ws2.Columns.Item(13).ColumnWidth = 16.87474886576337
$ws2.Columns.Item(14).ColumnWidth = 39.166666666666664

# ======================= de-de =======================
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Add($ws3.Cells.Item(5, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b836924ea5537688c547a030c04ea832cee94f00/e2e/e91463e1-37c1-4132-936d-4741263aa477.md", "", "", "e91463e1-37c1-4132-936d-4741263aa477.md")
$ws3.Cells.Item(5, 2).Value = ".md"
$ws3.Cells.Item(5, 3).Value = "Ready for handoff"
$ws3.Cells.Item(5, 4).Value = "e2e"
$ws3.Cells.Item(5, 5).Value = "ht"
$ws3.Cells.Item(5, 6).Value = "False"
$ws3.Cells.Item(5, 7).Value = "e91463e1-37c1-4132-936d-4741263aa477.eb0d211ad163a74bdec9c4fa15f72fa75ab483ad.de-de.xlf"
$ws3.Cells.Item(5, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(5, 8).Value = "2016-10-19 12:32:14"
$ws3.Cells.Item(5, 9).Value = ""
$ws3.Cells.Item(5, 10).Value = ""
$ws3.Cells.Item(5, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(5, 11).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(5, 12).Value = ""
$ws3.Cells.Item(5, 13).Value = "True"
$ws3.Cells.Item(5, 14).Value = ""
$ws3.Cells.Item(5, 15).Value = "False"
$ws3.Cells.Item(5, 16).Value = ""

$ws3.Hyperlinks.Add($ws3.Cells.Item(6, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b836924ea5537688c547a030c04ea832cee94f00/e2e/72a22f36-255e-4335-94e1-dc1ccefda6fd.png", "", "", "72a22f36-255e-4335-94e1-dc1ccefda6fd.png")
$ws3.Cells.Item(6, 2).Value = ".png"
$ws3.Cells.Item(6, 3).Value = "Ready for handoff"
$ws3.Cells.Item(6, 4).Value = "e2e"
$ws3.Cells.Item(6, 5).Value = "ht"
$ws3.Cells.Item(6, 6).Value = "False"
$ws3.Cells.Item(6, 7).Value = "e80e7f845ffd185e6ca385a6da93685d70922b36.png"
$ws3.Cells.Item(6, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(6, 8).Value = "2016-10-19 12:32:14"
$ws3.Cells.Item(6, 9).Value = ""
$ws3.Cells.Item(6, 10).Value = ""
$ws3.Cells.Item(6, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(6, 11).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(6, 12).Value = ""
$ws3.Cells.Item(6, 13).Value = "True(Dependency)"
$ws3.Cells.Item(6, 14).Value = "e2e\e91463e1-37c1-4132-936d-4741263aa477.md"
$ws3.Cells.Item(6, 15).Value = "False"
$ws3.Cells.Item(6, 16).Value = ""

$ws3.Hyperlinks.Add($ws3.Cells.Item(7, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b836924ea5537688c547a030c04ea832cee94f00/e2e/86614a44-b657-4e47-b440-569090e5a0bf.png", "", "", "86614a44-b657-4e47-b440-569090e5a0bf.png")
$ws3.Cells.Item(7, 2).Value = ".png"
$ws3.Cells.Item(7, 3).Value = "Ready for handoff"
$ws3.Cells.Item(7, 4).Value = "e2e"
$ws3.Cells.Item(7, 5).Value = "ht"
$ws3.Cells.Item(7, 6).Value = "False"
$ws3.Cells.Item(7, 7).Value = "3c4c472bb964b0033e821dd8fb8d4d40e13f94cd.png"
$ws3.Cells.Item(7, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(7, 8).Value = "2016-10-19 12:32:14"
$ws3.Cells.Item(7, 9).Value = ""
$ws3.Cells.Item(7, 10).Value = ""
$ws3.Cells.Item(7, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(7, 11).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(7, 12).Value = ""
$ws3.Cells.Item(7, 13).Value = "True(Dependency)"
$ws3.Cells.Item(7, 14).Value = "e2e\e91463e1-37c1-4132-936d-4741263aa477.md"
$ws3.Cells.Item(7, 15).Value = "False"
$ws3.Cells.Item(7, 16).Value = ""

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P7"))
$ws3.Columns.Item(13).ColumnWidth = 16.87474886576337
$ws3.Columns.Item(14).ColumnWidth = 39.166666666666664
